$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update input values (formulas in column F will recalculate automatically)
$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 16.5
$ws.Range("B5").Value = 11.3

# Update the active selection to B1 (as recorded in the saved view state)
$ws.Range("B1").Select()
